$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.169212666666667"
$ws.Range("H2").Value = [double]"3.507638"
$ws.Range("I2").Value = [double]"0.0005193657195729173"
$ws.Range("J2").Value = [double]"0.0005193657195729173"
$ws.Range("M2").Value = [double]"0.141694"
$ws.Range("N2").Value = [double]"0.425082"
$ws.Range("O2").Value = [double]"0.01763793963212447"
$ws.Range("P2").Value = [double]"0.01763793963212447"
$ws.Range("Q2").Value = [double]"0.1656704195906666"
$ws.Range("R2").Value = [double]"1.491033776316"
$ws.Range("S2").Value = [double]"9.160541208822002E-06"
$ws.Range("T2").Value = [double]"9.160541208822001E-06"
# Row 3
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.169212666666667"
$ws.Range("H3").Value = [double]"3.507638"
$ws.Range("I3").Value = [double]"0.0005193657195729173"
$ws.Range("J3").Value = [double]"0.0005193657195729173"
$ws.Range("O3").Value = [double]"0.2714637835982539"
$ws.Range("P3").Value = [double]"0.2714637835982538"
$ws.Range("Q3").Value = [double]"2.549817034778889"
$ws.Range("R3").Value = [double]"22.94835331301"
$ws.Range("S3").Value = [double]"0.0001409889833064939"
$ws.Range("T3").Value = [double]"0.0001409889833064938"
# Row 4
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.169212666666667"
$ws.Range("H4").Value = [double]"3.507638"
$ws.Range("I4").Value = [double]"0.0005193657195729173"
$ws.Range("J4").Value = [double]"0.0005193657195729173"
$ws.Range("M4").Value = [double]"5.710985666666667"
$ws.Range("N4").Value = [double]"17.132957"
$ws.Range("O4").Value = [double]"0.7108982767696218"
$ws.Range("P4").Value = [double]"0.7108982767696217"
$ws.Range("Q4").Value = [double]"6.677356780618444"
$ws.Range("R4").Value = [double]"60.09621102556601"
$ws.Range("S4").Value = [double]"0.0003692161950576016"
$ws.Range("T4").Value = [double]"0.0003692161950576015"
# Row 5
$ws.Range("I5").Value = [double]"0.9638330474556795"
$ws.Range("J5").Value = [double]"0.9638330474556795"
$ws.Range("M5").Value = [double]"0.141694"
$ws.Range("N5").Value = [double]"0.425082"
$ws.Range("O5").Value = [double]"0.01763793963212447"
$ws.Range("P5").Value = [double]"0.01763793963212447"
$ws.Range("Q5").Value = [double]"307.449297036084"
$ws.Range("R5").Value = [double]"2767.043673324756"
$ws.Range("S5").Value = [double]"0.01700002910646984"
$ws.Range("T5").Value = [double]"0.01700002910646983"
# Row 6
$ws.Range("I6").Value = [double]"0.9638330474556795"
$ws.Range("J6").Value = [double]"0.9638330474556795"
$ws.Range("O6").Value = [double]"0.2714637835982539"
$ws.Range("P6").Value = [double]"0.2714637835982538"
$ws.Range("S6").Value = [double]"0.2616457658193542"
$ws.Range("T6").Value = [double]"0.2616457658193541"
# Row 7
$ws.Range("I7").Value = [double]"0.9638330474556795"
$ws.Range("J7").Value = [double]"0.9638330474556795"
$ws.Range("M7").Value = [double]"5.710985666666667"
$ws.Range("N7").Value = [double]"17.132957"
$ws.Range("O7").Value = [double]"0.7108982767696218"
$ws.Range("P7").Value = [double]"0.7108982767696217"
$ws.Range("Q7").Value = [double]"12391.76343811184"
$ws.Range("R7").Value = [double]"111525.8709430065"
$ws.Range("S7").Value = [double]"0.6851872525298557"
$ws.Range("T7").Value = [double]"0.6851872525298556"
# Row 8
$ws.Range("G8").Value = [double]"80.250984"
$ws.Range("H8").Value = [double]"240.752952"
$ws.Range("I8").Value = [double]"0.03564758682474761"
$ws.Range("J8").Value = [double]"0.0356475868247476"
$ws.Range("M8").Value = [double]"0.141694"
$ws.Range("N8").Value = [double]"0.425082"
$ws.Range("O8").Value = [double]"0.01763793963212447"
$ws.Range("P8").Value = [double]"0.01763793963212447"
$ws.Range("Q8").Value = [double]"11.371082926896"
$ws.Range("R8").Value = [double]"102.339746342064"
$ws.Range("S8").Value = [double]"0.0006287499844458139"
$ws.Range("T8").Value = [double]"0.0006287499844458136"
# Row 9
$ws.Range("G9").Value = [double]"80.250984"
$ws.Range("H9").Value = [double]"240.752952"
$ws.Range("I9").Value = [double]"0.03564758682474761"
$ws.Range("J9").Value = [double]"0.0356475868247476"
$ws.Range("O9").Value = [double]"0.2714637835982539"
$ws.Range("P9").Value = [double]"0.2714637835982538"
$ws.Range("Q9").Value = [double]"175.01121215556"
$ws.Range("R9").Value = [double]"1575.10090940004"
$ws.Range("S9").Value = [double]"0.009677028795593252"
$ws.Range("T9").Value = [double]"0.009677028795593247"
# Row 10
$ws.Range("G10").Value = [double]"80.250984"
$ws.Range("H10").Value = [double]"240.752952"
$ws.Range("I10").Value = [double]"0.03564758682474761"
$ws.Range("J10").Value = [double]"0.0356475868247476"
$ws.Range("M10").Value = [double]"5.710985666666667"
$ws.Range("N10").Value = [double]"17.132957"
$ws.Range("O10").Value = [double]"0.7108982767696218"
$ws.Range("P10").Value = [double]"0.7108982767696217"
$ws.Range("Q10").Value = [double]"458.312219359896"
$ws.Range("R10").Value = [double]"4124.809974239064"
$ws.Range("S10").Value = [double]"0.02534180804470855"
$ws.Range("T10").Value = [double]"0.02534180804470854"
